$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency market data scraped on 2024-11-03.
# Numeric-looking text values are prefixed with a literal apostrophe so that
# Excel keeps them as text (matching the workbook's original inline-string cells)
# instead of auto-converting them to numbers.

$ws.Range("D2").Value = '''68.107.43'
$ws.Range("E2").Value = '''  -1.93%  '
$ws.Range("D3").Value = '''2.437.45'
$ws.Range("E3").Value = '''  -1.74%  '
$ws.Range("E4").Value = '''  -0.08%  '
$ws.Range("D5").Value = '''554.61'
$ws.Range("E5").Value = '''  -1.97%  '
$ws.Range("D6").Value = '''161.37'
$ws.Range("E6").Value = '''  -1.86%  '
$ws.Range("E7").Value = '''  -0.03%  '
$ws.Range("D8").Value = '''0.497'
$ws.Range("E8").Value = '''  -2.56%  '
$ws.Range("D9").Value = '''2.435.40'
$ws.Range("E10").Value = '''  -6.18%  '
$ws.Range("E11").Value = '''  -1.95%  '
$ws.Range("E12").Value = '''  -5.45%  '
$ws.Range("D13").Value = '''4.71'
$ws.Range("E13").Value = '''  -3.64%  '
$ws.Range("D14").Value = '''2.883.99'
$ws.Range("E14").Value = '''  -1.81%  '
$ws.Range("D15").Value = '''68.001.72'
$ws.Range("E15").Value = '''  -2.08%  '
$ws.Range("E16").Value = '''  -4.69%  '
$ws.Range("D17").Value = '''22.92'
$ws.Range("E17").Value = '''  -5.02%  '
$ws.Range("D18").Value = '''2.429.35'
$ws.Range("E18").Value = '''  -2.38%  '
$ws.Range("D19").Value = '''10.74'
$ws.Range("E19").Value = '''  -3.63%  '
$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").Value = '''336.39'
$ws.Range("E20").Value = '''  -2.30%  '
$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").Value = '''7.02'
$ws.Range("E21").Value = '''  -4.17%  '
$ws.Range("E22").Value = '''  -3.47%  '
$ws.Range("E23").Value = '''  -0.27%  '
$ws.Range("D24").Value = '''1.81'
$ws.Range("E24").Value = '''  -5.02%  '
$ws.Range("D25").Value = '''66.72'
$ws.Range("E25").Value = '''  -4.79%  '
$ws.Range("D26").Value = '''2.562.95'
$ws.Range("E26").Value = '''  -1.78%  '
$ws.Range("E27").Value = '''  -6.92%  '
$ws.Range("D28").Value = '''1.00'
$ws.Range("E28").Value = '''  +0.18%  '
$ws.Range("D29").Value = '''7.97'
$ws.Range("E29").Value = '''  -7.31%  '
$ws.Range("E30").Value = '''  -6.10%  '
$ws.Range("D31").Value = '''7.03'
$ws.Range("E31").Value = '''  -8.67%  '
$ws.Range("D32").Value = '''0.999'
$ws.Range("E32").Value = '''  -0.09%  '
$ws.Range("D33").Value = '''418.91'
$ws.Range("E33").Value = '''  -4.75%  '
$ws.Range("E34").Value = '''  -4.46%  '
$ws.Range("D35").Value = '''1.61'
$ws.Range("E35").Value = '''  -5.10%  '
$ws.Range("D36").Value = '''157.73'
$ws.Range("E36").Value = '''  +0.82%  '
$ws.Range("D37").Value = '''18.97'
$ws.Range("E37").Value = '''  -0.37%  '
$ws.Range("E38").Value = '''  +0.02%  '
$ws.Range("D39").Value = '''0.107'
$ws.Range("E39").Value = '''  -4.66%  '
$ws.Range("D40").Value = '''17.57'
$ws.Range("E40").Value = '''  -2.83%  '
$ws.Range("D41").Value = '''0.298'
$ws.Range("E41").Value = '''  -4.44%  '
$ws.Range("E42").Value = '''  -4.75%  '
$ws.Range("E43").Value = '''  -6.48%  '
$ws.Range("E44").Value = '''  +0.21%  '
$ws.Range("E45").Value = '''  -5.75%  '
$ws.Range("D46").Value = '''132.64'
$ws.Range("E46").Value = '''  -4.59%  '
$ws.Range("D47").Value = '''3.28'
$ws.Range("E47").Value = '''  -4.03%  '
$ws.Range("E48").Value = '''  -2.33%  '
$ws.Range("E49").Value = '''  -7.42%  '
$ws.Range("E50").Value = '''  -2.78%  '
$ws.Range("D51").Value = '''0.0898'
$ws.Range("E51").Value = '''  -2.09%  '

Write-Output "Updated $(84) cells with refreshed crypto market data."
